$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header subject numbers) tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) tweaks
$ws.Range("B2").Value = 1.2848039557918909
$ws.Range("C2").Value = 1.4438791746789432
$ws.Range("D2").Value = 3.854618248298646
$ws.Range("E2").Value = 1.095336699717498

# Row 3 (STR) tweaks
$ws.Range("B3").Value = 2.0886367795514373
$ws.Range("C3").Value = 0.83904061698529964
$ws.Range("D3").Value = 3.3165995798502315
$ws.Range("E3").Value = 0.30385789479096725

# Update the active selection to match the edited range
$ws.Range("B1:E3").Select()
